$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the "Absent" (H) column flips from 0 to 1
$absentRows = @(3, 6, 7, 8, 9, 13, 14, 15, 16, 17, 18)
foreach ($r in $absentRows) {
    $ws.Range("H$r").Value = 1
}

# Rows where "Total Attendance Count" (D) and "Real" (E) flip from 0 to 1
$presentRows = @(4, 5, 10, 11, 12)
foreach ($r in $presentRows) {
    $ws.Range("D$r").Value = 1
    $ws.Range("E$r").Value = 1
}
